$wb = $excel.ActiveWorkbook

# --- Add the new "New Customer" worksheet after "Email" -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "New Customer"

# --- Header row -------------------------------------------------------------
$newSheet.Range("A1").Value = "Customer Name"
$newSheet.Range("B1").Value = "Gender"

$newSheet.Columns(3).NumberFormat = "mm-dd-yy"
$newSheet.Range("C1").ClearFormats()
$newSheet.Range("C1").Value = "Date Of Birth"

$newSheet.Range("D1").Value = "Address"
$newSheet.Range("E1").Value = "City "
$newSheet.Range("F1").Value = "State"
$newSheet.Range("G1").Value = "PIN"
$newSheet.Range("H1").Value = "Tel. No"
$newSheet.Range("I1").Value = "Email"

# --- Row 2: Rajesh Kumar -----------------------------------------------------
$newSheet.Range("A2").Value = "Rajesh Kumar"
$newSheet.Range("B2").Value = "Male"
$newSheet.Range("C2").Value = [DateTime]"1991-10-20"
$newSheet.Range("D2").Value = "Hno-104,Street no.2"
$newSheet.Range("E2").Value = "Mumbai"
$newSheet.Range("F2").Value = "Maharashtra"
$newSheet.Range("G2").Value = 400001
$newSheet.Range("H2").Value = 23452472348
$newSheet.Range("I2").Value = "rkumar12@sdhaj.com"

# --- Row 3: Saniya Patel -----------------------------------------------------
$newSheet.Range("A3").Value = "Saniya Patel"
$newSheet.Range("B3").Value = "Female"
$newSheet.Range("C3").Value = [DateTime]"1993-09-15"
$newSheet.Range("D3").Value = "Fl.A-405, Lane2"
$newSheet.Range("E3").Value = "Pune"
$newSheet.Range("F3").Value = "Maharashtra"
$newSheet.Range("G3").Value = 416532
$newSheet.Range("H3").Value = 234587912
$newSheet.Range("I3").Value = "sanpatel@asa.com"

# --- Column widths (auto-fit to content, like Excel would on data entry) ---
$newSheet.Columns("A").AutoFit()
$newSheet.Columns("B").AutoFit()
$newSheet.Columns("C").AutoFit()
$newSheet.Columns("H").AutoFit()

# --- Blank placeholder row added to the Login sheet -------------------------
$ws1 = $wb.Worksheets.Item("Login")
$ws1.Range("A2").Value = "  "
$ws1.Range("B2").Value = "  "
$ws1.Range("A2").Select()

# --- Final selections / active sheet ----------------------------------------
$newSheet.Range("H27").Select()

$ws2 = $wb.Worksheets.Item("Email")
$ws2.Range("C20").Select()
